# Apply the "Use US EPS values" note update to the RTMF workbook's About sheet.
#
# Summary of the change (per the commit message / diff):
#   - A new note "The EU EPS uses values from the US EPS." is inserted on the
#     "About" sheet, just below the "Source:" / notes block (between the old
#     row 13 and the paragraph that used to start at row 15), pushing the
#     following "non-motorized modes" paragraph and the rest of the sheet
#     down by two rows.
#   - The previously selected cells on the "About" and "RTMF-freight" sheets
#     are updated to reflect where the editing author's cursor ended up.

$wb = $excel.ActiveWorkbook

$about = $wb.Worksheets.Item("About")

# Insert two new blank rows at row 15 (shifts the old rows 15-28 down to 17-30).
$about.Range("A15:A16").EntireRow.Insert()

# Row 15 gets the new note; row 16 stays blank, matching the gap pattern
# already used elsewhere on this sheet (e.g. the old blank row 14/20).
$about.Range("A15").Value = "The EU EPS uses values from the US EPS."

# Restore the "About" sheet as the active sheet/selection.
[void]$about.Select()
[void]$about.Range("L24").Select()

# Update the remembered selection on the "RTMF-freight" sheet.
$freight = $wb.Worksheets.Item("RTMF-freight")
[void]$freight.Range("B30").Select()

# Leave "About" as the active sheet/tab when the workbook is saved.
[void]$about.Select()
